$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the numeric-looking Price column cells being updated,
# so Excel does not auto-convert them to numbers (source data is stored as text).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "29.418.01"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.876.01"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "0.7130"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").Value = "241.61"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.07842"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").Value = "25.20"
$ws.Range("E10").Value = "  +7.52%  "
$ws.Range("D11").Value = "0.08248"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "0.7286"
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.876.41"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.258"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "90.85"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "29.418.21"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "5.902"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "247.00"
$ws.Range("E18").Value = "  +4.06%  "
$ws.Range("D19").Value = "0.000007864"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "7.942"
$ws.Range("E22").Value = "  +6.74%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "0.1577"
$ws.Range("E24").Value = "  +10.44%  "
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "8.987"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").Value = "1.362"
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("D30").Value = "4.359"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "4.124"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").Value = "0.05308"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").Value = "1.925"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "1.198"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").Value = "0.7212"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "2.680"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "0.01860"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "1.255.42"
$ws.Range("E38").Value = "  +9.30%  "
$ws.Range("D39").Value = "2.728"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "0.9083"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "73.88"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("D42").Value = "6.139"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "103.27"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "2.010.98"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").Value = "2.932"
$ws.Range("E47").Value = "  +13.55%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "0.4314"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "9.241"
$ws.Range("E51").Value = "  +1.07%  "

# Restore default (Normal) style on the cells we forced to text format,
# so no extra explicit style index lingers on them.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
